$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.226.83"
$ws.Range("E2").Value = "'  -0.16%  "
$ws.Range("D3").Value = "'1.795.27"
$ws.Range("E3").Value = "'  +0.20%  "
$ws.Range("E4").Value = "'  +0.25%  "
$ws.Range("D5").Value = "'226.54"
$ws.Range("E5").Value = "'  -0.04%  "
$ws.Range("D6").Value = "'0.569"
$ws.Range("E6").Value = "'  +2.61%  "
$ws.Range("E7").Value = "'  +0.16%  "
$ws.Range("D8").Value = "'35.81"
$ws.Range("E8").Value = "'  +9.78%  "
$ws.Range("D9").Value = "'0.297"
$ws.Range("E9").Value = "'  +0.93%  "
$ws.Range("D10").Value = "'0.0686"
$ws.Range("E10").Value = "'  -0.30%  "
$ws.Range("D11").Value = "'0.0962"
$ws.Range("E11").Value = "'  +1.71%  "
$ws.Range("D12").Value = "'2.056.90"
$ws.Range("E12").Value = "'  +0.31%  "
$ws.Range("D13").Value = "'11.40"
$ws.Range("E13").Value = "'  +3.67%  "
$ws.Range("D14").Value = "'1.791.14"
$ws.Range("E14").Value = "'  -0.07%  "
$ws.Range("D15").Value = "'0.637"
$ws.Range("E15").Value = "'  +1.06%  "
$ws.Range("D16").Value = "'4.43"
$ws.Range("E16").Value = "'  +4.08%  "
$ws.Range("D17").Value = "'34.244.87"
$ws.Range("E17").Value = "'  -0.17%  "
$ws.Range("D18").Value = "'68.62"
$ws.Range("E18").Value = "'  +0.48%  "
$ws.Range("D19").Value = "'243.31"
$ws.Range("E19").Value = "'  -0.07%  "
$ws.Range("D20").Value = "'0.0₃0786"
$ws.Range("E20").Value = "'  -0.74%  "
$ws.Range("D21").Value = "'11.45"
$ws.Range("E21").Value = "'  +2.22%  "
$ws.Range("E22").Value = "'  +0.13%  "
$ws.Range("D23").Value = "'4.13"
$ws.Range("E23").Value = "'  -0.23%  "
$ws.Range("D24").Value = "'2.14"
$ws.Range("E24").Value = "'  +3.63%  "
$ws.Range("D25").Value = "'171.54"
$ws.Range("E25").Value = "'  +3.56%  "
$ws.Range("D26").Value = "'7.82"
$ws.Range("E26").Value = "'  +7.51%  "
$ws.Range("D27").Value = "'16.66"
$ws.Range("E27").Value = "'  +1.39%  "
$ws.Range("E28").Value = "'  +1.40%  "
$ws.Range("E29").Value = "'  +0.06%  "
$ws.Range("D30").Value = "'3.96"
$ws.Range("E30").Value = "'  -0.17%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.24"
$ws.Range("E31").Value = "'  +0.37%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0525"
$ws.Range("E32").Value = "'  +0.43%  "
$ws.Range("D33").Value = "'3.79"
$ws.Range("E33").Value = "'  +0.42%  "
$ws.Range("E34").Value = "'  -0.54%  "
$ws.Range("D35").Value = "'1.389.07"
$ws.Range("E35").Value = "'  -0.84%  "
$ws.Range("D36").Value = "'0.666"
$ws.Range("E36").Value = "'  -0.05%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.06"
$ws.Range("E37").Value = "'  +0.03%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.40"
$ws.Range("E38").Value = "'  -6.64%  "
$ws.Range("E39").Value = "'  -0.22%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.82"
$ws.Range("E40").Value = "'  -0.35%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'81.83"
$ws.Range("E41").Value = "'  -3.06%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "'  +0.39%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.949"
$ws.Range("E43").Value = "'  +1.84%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.20"
$ws.Range("E44").Value = "'  +7.48%  "
$ws.Range("D45").Value = "'13.20"
$ws.Range("E45").Value = "'  -4.36%  "
$ws.Range("D46").Value = "'0.0504"
$ws.Range("E46").Value = "'  -3.64%  "
$ws.Range("D47").Value = "'6.00"
$ws.Range("E47").Value = "'  +0.40%  "
$ws.Range("D48").Value = "'1.957.49"
$ws.Range("E48").Value = "'  +0.37%  "
$ws.Range("E49").Value = "'  +0.09%  "
$ws.Range("D50").Value = "'103.57"
$ws.Range("E50").Value = "'  -0.88%  "
$ws.Range("D51").Value = "'0.0₆0126"
$ws.Range("E51").Value = "'  -1.54%  "
